# Update workbook values per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

# Sheet 1 = 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 9819
$ws1.Range("F7").Value = 350
$ws1.Range("F8").Value = 383
$ws1.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202404/UHohfLBe1714358642039.jpeg"
$ws1.Range("F14").Value = 12452
$ws1.Range("F24").Value = 2744
$ws1.Range("F25").Value = 2109
$ws1.Range("F29").Value = 2160
$ws1.Range("F30").Value = 1058
$ws1.Range("F31").Value = 4230
$ws1.Range("F32").Value = 3735
$ws1.Range("F33").Value = 741
$ws1.Range("F34").Value = 2640
$ws1.Range("F36").Value = 51
$ws1.Range("F37").Value = 1344
$ws1.Range("F39").Value = 782
$ws1.Range("F42").Value = 462
$ws1.Range("F43").Value = 603
$ws1.Range("F49").Value = 155

# Sheet 2 = 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F17").Value = 10

# Sheet 4 = 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 9819
$ws4.Range("F10").Value = 350
$ws4.Range("F11").Value = 383
$ws4.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202404/UHohfLBe1714358642039.jpeg"
$ws4.Range("F16").Value = 12452
$ws4.Range("F24").Value = 2744
$ws4.Range("F25").Value = 2109
$ws4.Range("F28").Value = 2160
$ws4.Range("F29").Value = 1058
$ws4.Range("F30").Value = 4230
$ws4.Range("F31").Value = 3735
$ws4.Range("F32").Value = 741
$ws4.Range("F33").Value = 2640
$ws4.Range("F35").Value = 51
$ws4.Range("F36").Value = 1344
$ws4.Range("F38").Value = 782
$ws4.Range("F41").Value = 462
$ws4.Range("F43").Value = 603
$ws4.Range("F49").Value = 155
